$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear E column values for rows 3,4,5,6,8,9,10,11,13,14 (keep style)
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""

# Update D13 value
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
